$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.368.19'
$ws.Range('E2').Value = '  +3.47%  '
$ws.Range('D3').Value = '2.257.08'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'296.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = "'87.98"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.77%  '
$ws.Range('D7').Value = "'0.517"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.85%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.478"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('D10').Value = "'31.53"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.24%  '
$ws.Range('D11').Value = "'0.0805"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.49%  '
$ws.Range('D12').Value = "'47.35"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = "'6.50"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.72%  '
$ws.Range('D15').Value = '2.595.40'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = "'14.38"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').Value = '2.237.22'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = "'0.739"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('D19').Value = '40.229.37'
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').Value = '0.0₃0898'
$ws.Range('E20').Value = '  +4.52%  '
$ws.Range('D21').Value = "'5.90"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('D22').Value = "'10.78"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.30%  '
$ws.Range('D23').Value = "'66.09"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').Value = "'237.68"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.85%  '
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').Value = "'2.49"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('D27').Value = "'1.86"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.64%  '
$ws.Range('D28').Value = "'23.36"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.20%  '
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('D30').Value = "'9.33"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.74%  '
$ws.Range('D31').Value = "'33.98"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.27%  '
$ws.Range('D32').Value = "'153.43"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'0.999"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = "'4.95"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.86%  '
$ws.Range('D35').Value = "'0.0723"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.86%  '
$ws.Range('D36').Value = "'2.40"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('D37').Value = "'16.87"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +15.81%  '
$ws.Range('D38').Value = "'0.103"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.86%  '
$ws.Range('E39').Value = '  +3.07%  '
$ws.Range('D40').Value = "'2.75"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.57%  '
$ws.Range('D41').Value = "'1.71"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.03%  '
$ws.Range('D42').Value = "'3.85"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.85%  '
$ws.Range('D43').Value = '2.023.59'
$ws.Range('E43').Value = '  +6.47%  '
$ws.Range('E44').Value = '  +7.13%  '
$ws.Range('D45').Value = "'0.0273"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.67%  '
$ws.Range('D46').Value = "'10.02"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.44%  '
$ws.Range('D47').Value = "'16.68"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('D48').Value = "'2.62"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.25%  '
$ws.Range('D49').Value = '2.466.54'
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('D50').Value = "'72.00"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.20%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = "'1.47"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.85%  '
